$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Project Total Costs")
$ws.Range("B2").Value = 881263.5856249999
$ws.Range("B3").Value = 181390.125184
$ws.Range("B4").Value = 708378.293358
$ws.Range("B5").Value = 699873.4604409999
$ws.Range("B7").Value = 0.6805021261686089

$ws = $wb.Worksheets.Item("Components Capacity and Cost")
$ws.Range("B2").Value = 533.086236485
$ws.Range("B3").Value = 28.0283032253
$ws.Range("B4").Value = 0.016150133392
$ws.Range("B5").Value = 171.934662575
$ws.Range("B6").Value = 20.8375089279
$ws.Range("B7").Value = 287866.5677019
$ws.Range("B8").Value = 5605.66064506
$ws.Range("B9").Value = 19.3801600704
$ws.Range("B10").Value = 343869.32515
$ws.Range("B11").Value = 62512.5267837
$ws.Range("B12").Value = 699873.4604407303

$ws = $wb.Worksheets.Item("Yearly Costs Info")
$ws.Range("B2").Value = 5757.331354038
$ws.Range("C2").Value = 569.2871365376801
$ws.Range("D2").Value = 8127.637038674
$ws.Range("E2").Value = 14454.25552924968
$ws.Range("F2").Value = 10.26666487699951
$ws.Range("G2").Value = 8901.973929389767
$ws.Range("B3").Value = 5757.331354038
$ws.Range("C3").Value = 569.2871365376801
$ws.Range("D3").Value = 8127.637038674
$ws.Range("E3").Value = 14454.25552924968
$ws.Range("F3").Value = 8.713501243434349
$ws.Range("G3").Value = 8921.29360126965
$ws.Range("B4").Value = 5757.331354038
$ws.Range("C4").Value = 569.2871365376801
$ws.Range("D4").Value = 8127.637038674
$ws.Range("E4").Value = 14454.25552924968
$ws.Range("F4").Value = 72.62322237708076
$ws.Range("G4").Value = 11037.6794280729
$ws.Range("B5").Value = 5757.331354038
$ws.Range("C5").Value = 569.2871365376801
$ws.Range("D5").Value = 8127.637038674
$ws.Range("E5").Value = 14454.25552924968
$ws.Range("F5").Value = 91.01999848860085
$ws.Range("G5").Value = 11035.3890742394
$ws.Range("B6").Value = 5757.331354038
$ws.Range("C6").Value = 569.2871365376801
$ws.Range("D6").Value = 8127.637038674
$ws.Range("E6").Value = 14454.25552924968
$ws.Range("F6").Value = 90.89915593917038
$ws.Range("G6").Value = 11035.41157114422
$ws.Range("B7").Value = 5757.331354038
$ws.Range("C7").Value = 569.2871365376801
$ws.Range("D7").Value = 8127.637038674
$ws.Range("E7").Value = 14454.25552924968
$ws.Range("F7").Value = 92.29260687403114
$ws.Range("G7").Value = 11035.15215792765
$ws.Range("B8").Value = 5757.331354038
$ws.Range("C8").Value = 569.2871365376801
$ws.Range("D8").Value = 8127.637038674
$ws.Range("E8").Value = 14454.25552924968
$ws.Range("F8").Value = 86.06974713480078
$ws.Range("G8").Value = 11036.31439383156
$ws.Range("B9").Value = 5757.331354038
$ws.Range("C9").Value = 569.2871365376801
$ws.Range("D9").Value = 8127.637038674
$ws.Range("E9").Value = 14454.25552924968
$ws.Range("F9").Value = 3766.661691524581
$ws.Range("G9").Value = 14754.89891929363
$ws.Range("B10").Value = 5757.331354038
$ws.Range("C10").Value = 569.2871365376801
$ws.Range("D10").Value = 8127.637038674
$ws.Range("E10").Value = 14454.25552924968
$ws.Range("F10").Value = 3732.57922245498
$ws.Range("G10").Value = 14761.24443967141
$ws.Range("B11").Value = 5757.331354038
$ws.Range("C11").Value = 569.2871365376801
$ws.Range("D11").Value = 8127.637038674
$ws.Range("E11").Value = 14454.25552924968
$ws.Range("F11").Value = 3739.642232258679
$ws.Range("G11").Value = 14759.92952075165
$ws.Range("B12").Value = 5757.331354038
$ws.Range("C12").Value = 569.2871365376801
$ws.Range("D12").Value = 8127.637038674
$ws.Range("E12").Value = 14454.25552924968
$ws.Range("F12").Value = 3739.627197121625
$ws.Range("G12").Value = 14759.93231987801
$ws.Range("B13").Value = 5757.331354038
$ws.Range("C13").Value = 569.2871365376801
$ws.Range("D13").Value = 8127.637038674
$ws.Range("E13").Value = 14454.25552924968
$ws.Range("F13").Value = 3739.629627440963
$ws.Range("G13").Value = 14759.93186743954
$ws.Range("B14").Value = 5757.331354038
$ws.Range("C14").Value = 569.2871365376801
$ws.Range("D14").Value = 8127.637038674
$ws.Range("E14").Value = 14454.25552924968
$ws.Range("F14").Value = 3740.026234171446
$ws.Range("G14").Value = 14759.85803017139
$ws.Range("B15").Value = 5757.331354038
$ws.Range("C15").Value = 569.2871365376801
$ws.Range("D15").Value = 8127.637038674
$ws.Range("E15").Value = 14454.25552924968
$ws.Range("F15").Value = 3755.748012861778
$ws.Range("G15").Value = 14756.93226478665
$ws.Range("B16").Value = 5757.331354038
$ws.Range("C16").Value = 569.2871365376801
$ws.Range("D16").Value = 8127.637038674
$ws.Range("E16").Value = 14454.25552924968
$ws.Range("F16").Value = 21981.41623442561
$ws.Range("G16").Value = 14041.06698898695
$ws.Range("B17").Value = 5757.331354038
$ws.Range("C17").Value = 569.2871365376801
$ws.Range("D17").Value = 8127.637038674
$ws.Range("E17").Value = 14454.25552924968
$ws.Range("F17").Value = 22015.05464370237
$ws.Range("G17").Value = 14045.13175366127
$ws.Range("B18").Value = 5757.331354038
$ws.Range("C18").Value = 569.2871365376801
$ws.Range("D18").Value = 8127.637038674
$ws.Range("E18").Value = 14454.25552924968
$ws.Range("F18").Value = 22015.3401029741
$ws.Range("G18").Value = 14045.18489499017
$ws.Range("B19").Value = 5757.331354038
$ws.Range("C19").Value = 569.2871365376801
$ws.Range("D19").Value = 8127.637038674
$ws.Range("E19").Value = 14454.25552924968
$ws.Range("F19").Value = 22015.31552407677
$ws.Range("G19").Value = 14045.1803190899
$ws.Range("B20").Value = 5757.331354038
$ws.Range("C20").Value = 569.2871365376801
$ws.Range("D20").Value = 8127.637038674
$ws.Range("E20").Value = 14454.25552924968
$ws.Range("F20").Value = 22015.94663418807
$ws.Range("G20").Value = 14045.29781027312
$ws.Range("B21").Value = 5757.331354038
$ws.Range("C21").Value = 569.2871365376801
$ws.Range("D21").Value = 8127.637038674
$ws.Range("E21").Value = 14454.25552924968
$ws.Range("F21").Value = 22003.46475807035
$ws.Range("G21").Value = 14043.04247608964

$ws = $wb.Worksheets.Item("Yearly Energy Averages")
$ws.Range("B2").Value = 0.006149356996374144
$ws.Range("C2").Value = 0.02436318937933031
$ws.Range("D2").Value = 106.0371181402666
$ws.Range("E2").Value = 55.23760786637101
$ws.Range("B3").Value = 0.004590476351806553
$ws.Range("C3").Value = 0.02143693627413258
$ws.Range("D3").Value = 106.26722311305
$ws.Range("E3").Value = 55.14154327188311
$ws.Range("B4").Value = 0.1310368413605871
$ws.Range("C4").Value = 0.01848223917845942
$ws.Range("D4").Value = 105.9224020357453
$ws.Range("E4").Value = 43.09770419000218
$ws.Range("B5").Value = 0.1680828095592418
$ws.Range("C5").Value = 0.01851018965716523
$ws.Range("D5").Value = 105.9347456369335
$ws.Range("E5").Value = 43.0825140802494
$ws.Range("B6").Value = 0.1678393097182145
$ws.Range("C6").Value = 0.01851019538659774
$ws.Range("D6").Value = 105.9347456381713
$ws.Range("E6").Value = 43.08257034390343
$ws.Range("B7").Value = 0.1706469812070397
$ws.Range("C7").Value = 0.01851031828857077
$ws.Range("D7").Value = 105.9347456798529
$ws.Range("E7").Value = 43.08192153188168
$ws.Range("B8").Value = 0.1580776571893881
$ws.Range("C8").Value = 0.01854704447014374
$ws.Range("D8").Value = 105.9347620864378
$ws.Range("E8").Value = 43.08480875732312
$ws.Range("B9").Value = 5.156537758480843
$ws.Range("C9").Value = 0.02334000382108572
$ws.Range("D9").Value = 100.4404385496902
$ws.Range("E9").Value = 20.10651668097542
$ws.Range("B10").Value = 5.109701348608501
$ws.Range("C10").Value = 0.02334352120887271
$ws.Range("D10").Value = 100.4404306440267
$ws.Range("E10").Value = 20.11364800307939
$ws.Range("B11").Value = 5.119406772722357
$ws.Range("C11").Value = 0.0233435439888775
$ws.Range("D11").Value = 100.4404306355573
$ws.Range("E11").Value = 20.11217096947526
$ws.Range("B12").Value = 5.119386113538319
$ws.Range("C12").Value = 0.02334354285202083
$ws.Range("D12").Value = 100.4404306347396
$ws.Range("E12").Value = 20.11217411408981
$ws.Range("B13").Value = 5.11938945317551
$ws.Range("C13").Value = 0.02334354275544633
$ws.Range("D13").Value = 100.4404306346396
$ws.Range("E13").Value = 20.11217360596892
$ws.Range("B14").Value = 5.119934441563851
$ws.Range("C14").Value = 0.02334354030485455
$ws.Range("D14").Value = 100.4404306329202
$ws.Range("E14").Value = 20.11209067406496
$ws.Range("B15").Value = 5.141536360845534
$ws.Range("C15").Value = 0.02334563728929531
$ws.Range("D15").Value = 100.4404212188214
$ws.Range("E15").Value = 20.10881075678516
$ws.Range("B16").Value = 22.1339819558098
$ws.Range("C16").Value = 0.0229770276501979
$ws.Range("D16").Value = 81.70533706558732
$ws.Range("E16").Value = 9.933065620749911
$ws.Range("B17").Value = 22.16785422371817
$ws.Range("C17").Value = 0.02301172379000663
$ws.Range("D17").Value = 81.70534778390657
$ws.Range("E17").Value = 9.930136536273798
$ws.Range("B18").Value = 22.16814191306015
$ws.Range("C18").Value = 0.02301172158868695
$ws.Range("D18").Value = 81.70534777355519
$ws.Range("E18").Value = 9.93011177933306
$ws.Range("B19").Value = 22.1681171419177
$ws.Range("C19").Value = 0.02301172205523605
$ws.Range("D19").Value = 81.70534777424575
$ws.Range("E19").Value = 9.930113911144788
$ws.Range("B20").Value = 22.16875332052011
$ws.Range("C20").Value = 0.02301155031108055
$ws.Range("D20").Value = 81.70534761631041
$ws.Range("E20").Value = 9.930059297014294
$ws.Range("B21").Value = 22.15501801105849
$ws.Range("C21").Value = 0.02440822884000263
$ws.Range("D21").Value = 81.70554154032622
$ws.Range("E21").Value = 9.930953449960974
